$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data as described by the diff (rows 2-51)
$ws.Range("D2").Value = "'27.118.12"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "'1.558.87"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'209.78"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'0.490"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'21.98"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").Value = "'0.0870"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "'1.780.88"
$ws.Range("E12").Value = "  -0.24%  "
$ws.Range("D13").Value = "'1.545.24"
$ws.Range("E13").Value = "  -1.11%  "
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "'27.091.10"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'61.69"
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "'7.44"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0₃0700"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'215.79"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'9.19"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("D25").Value = "'152.91"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'15.00"
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").Value = "'0.106"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "'1.433.04"
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  +3.32%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -0.55%  "
$ws.Range("E38").Value = "  +0.34%  "
$ws.Range("E39").Value = "  -0.73%  "
$ws.Range("E40").Value = "  +1.36%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "'0.997"
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("D45").Value = "'64.05"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'1.73"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("D47").Value = "'1.693.02"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'85.37"
$ws.Range("E48").Value = "  -2.05%  "
$ws.Range("D49").Value = "'0.0524"
$ws.Range("E49").Value = "  +0.58%  "
$ws.Range("D50").Value = "'0.0₇0984"
$ws.Range("E50").Value = "  -1.56%  "
$ws.Range("D51").Value = "'0.0948"
$ws.Range("E51").Value = "  -0.97%  "

# Re-apply the default (unstyled) style to the Price cells so the quote-prefix
# flag used above doesn't leave a stray style on the cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
